# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders the "Periodo Mora" (column E) to ascending chronological order
# and updates the corresponding "Valor Mora" (column F) and
# "Salario Basico" (column G) figures for rows 16-57 on Hoja1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, period (text), valor mora, salario basico
$data = @(
    @(16, "1610", 14000, 781242),
    @(17, "1611", 14000, 781242),
    @(18, "1612", 14000, 781242),
    @(19, "1701", 21000, 781242),
    @(20, "1702", 21000, 781242),
    @(21, "1703", 21000, 781242),
    @(22, "1704", 21000, 781242),
    @(23, "1705", 21000, 781242),
    @(24, "1706", 21000, 781242),
    @(25, "1707", 21000, 781242),
    @(26, "1708", 21000, 781242),
    @(27, "1709", 21000, 781242),
    @(28, "1710", 21000, 781242),
    @(29, "1711", 21000, 781242),
    @(30, "1712", 21000, 781242),
    @(31, "1801", 28000, 781242),
    @(32, "1802", 28000, 781242),
    @(33, "1803", 28000, 781242),
    @(34, "1804", 28000, 781242),
    @(35, "1805", 28000, 781242),
    @(36, "1806", 28000, 781242),
    @(37, "1807", 28000, 781242),
    @(38, "1808", 28000, 781242),
    @(39, "1809", 31249, 781242),
    @(40, "1810", 31249, 781242),
    @(41, "1811", 31249, 781242),
    @(42, "1812", 31249, 781242),
    @(43, "1901", 31249, 781242),
    @(44, "1902", 31249, 781242),
    @(45, "1903", 31249, 781242),
    @(46, "1904", 31249, 781242),
    @(47, "1905", 31249, 781242),
    @(48, "1906", 31249, 781242),
    @(49, "1907", 31249, 781242),
    @(50, "1908", 31249, 781242),
    @(51, "1909", 31249, 781242),
    @(52, "1910", 31249, 781242),
    @(53, "1911", 31249, 781242),
    @(54, "1912", 31249, 781242),
    @(55, "2001", 31249, 781242),
    @(56, "2002", 31249, 781242),
    @(57, "2003", 31249, 781242)
)

foreach ($item in $data) {
    $row = $item[0]
    $period = $item[1]
    $valorMora = $item[2]
    $salarioBasico = $item[3]

    $ws.Cells.Item($row, 5).Value = $period
    $ws.Cells.Item($row, 6).Value = $valorMora
    $ws.Cells.Item($row, 7).Value = $salarioBasico
}
